# Insert two new weekly data rows at the top of the Zanahoria data block
# (rows 230-231), shifting all the existing rows down by two positions.
# This mirrors a new week of price reporting being prepended to the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows before the current row 230, pushing the rest of the
# table (old rows 230-337) down to rows 232-339.
$ws.Range("A230:A231").EntireRow.Insert()

# New row 230: "Primera" quality record for 2022-10-05 (serial 44839)
$ws.Range("A230").Value = 7
$ws.Range("B230").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C230").Value = "Ñuble"
$ws.Range("D230").Value = 44839
$ws.Range("E230").Value = 16
$ws.Range("F230").Value = 100114013
$ws.Range("G230").Value = "Zanahoria"
$ws.Range("H230").Value = "Sin especificar"
$ws.Range("I230").Value = "Primera"
$ws.Range("J230").Value = 100
$ws.Range("K230").Value = 10000
$ws.Range("L230").Value = 10000
$ws.Range("M230").Value = 10000
$ws.Range("N230").Value = "`$/saco 20 kilos"
$ws.Range("O230").Value = "Región de Ñuble"
$ws.Range("P230").Value = 500
$ws.Range("Q230").Value = 20
$ws.Range("R230").Value = "Hortaliza"

# New row 231: "Segunda" quality record for the same date (serial 44839)
$ws.Range("A231").Value = 7
$ws.Range("B231").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C231").Value = "Ñuble"
$ws.Range("D231").Value = 44839
$ws.Range("E231").Value = 16
$ws.Range("F231").Value = 100114013
$ws.Range("G231").Value = "Zanahoria"
$ws.Range("H231").Value = "Sin especificar"
$ws.Range("I231").Value = "Segunda"
$ws.Range("J231").Value = 150
$ws.Range("K231").Value = 9000
$ws.Range("L231").Value = 9000
$ws.Range("M231").Value = 9000
$ws.Range("N231").Value = "`$/saco 20 kilos"
$ws.Range("O231").Value = "Región de Ñuble"
$ws.Range("P231").Value = 450
$ws.Range("Q231").Value = 20
$ws.Range("R231").Value = "Hortaliza"
